$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25, columns C-N (380 kV case)
$data = @(
    @(8.352934928600247, 6.874829620302777, 12.52404744410618, 43.06867727401279, 52.01908916136161, 20.1201661971888, 33.26187852372878, 10.37567679562269, 18.75880835312677, 10.14928694085233, 20.11147755274974, 20.08572880150762),
    @(8.33798645102128, 6.865644555026096, 12.54265619072867, 43.14347315040508, 52.08742958452556, 20.17643991615019, 33.33296213345718, 10.39733661833669, 18.48488214230177, 10.16776491525243, 20.01908539897646, 20.14993719565026),
    @(8.329886249732555, 6.860560160050184, 12.55536838191547, 43.19922430488381, 52.14378767728983, 20.21457913421916, 33.38367162493798, 10.41144210932979, 18.31800532289172, 10.17982591920729, 19.96525922771995, 20.19128665618568),
    @(8.326858791069732, 6.858629140813188, 12.5608726453083, 43.22440848078037, 52.17036290358508, 20.23102188462809, 33.40610849431925, 10.4173934506012, 18.25041322428799, 10.18492123893362, 19.94407039834718, 20.20862246292575),
    @(8.326372663227884, 6.858317043860901, 12.56180620320012, 43.2287390002342, 52.17499330104951, 20.23380655394692, 33.40994104300375, 10.41839395603495, 18.23921682427024, 10.18577821979454, 19.94059749625931, 20.21153043010494),
    @(8.329844310280336, 6.860533545429289, 12.55544130207572, 43.19955397445393, 52.14413148432839, 20.21479724168809, 33.38396704575889, 10.41152154763466, 18.31709197914228, 10.17989390558306, 19.96497042796993, 20.19151848499305),
    @(8.347558340657182, 6.871548147738296, 12.53019705034119, 43.09242421032068, 52.03965854976581, 20.13882416891787, 33.2849197808269, 10.38297805882588, 18.66413423154591, 10.15550992858003, 20.07902795793779, 20.10746914276),
    @(8.390743264154146, 6.897495878808243, 12.49087853303274, 42.96054014582279, 51.94945587229267, 20.01835677224314, 33.14690076035473, 10.33337936062529, 19.35156995451064, 10.11334926910457, 20.32499715660705, 19.95785819181177),
    @(8.427468624710968, 6.919131215007619, 12.46817185400305, 42.91159202811544, 51.95357923163974, 19.94730664814893, 33.07997273642889, 10.30079447672761, 19.85613060806785, 10.08579364164016, 20.51827727062733, 19.85711765214477),
    @(8.445225353662925, 6.929514179468383, 12.45917773853954, 42.89977880035131, 51.97080698454226, 19.91878957977297, 33.05704769468695, 10.28680134373195, 20.08459305355643, 10.07399445492474, 20.60870000051689, 19.81326144430098),
    @(8.452097081587473, 6.933522098345883, 12.45596336045072, 42.8968106978253, 51.97954030799312, 19.90853905282727, 33.04945031439795, 10.28162135090987, 20.17087487743938, 10.06963178910237, 20.64327930022524, 19.79693623346073),
    @(8.450610620993288, 6.93265556319864, 12.45664712567332, 42.89738295850075, 51.97756115119208, 19.91072228349585, 33.05103830966232, 10.28273167367452, 20.152304084635, 10.07056668464723, 20.63581733070043, 19.80043962952027),
    @(8.445787755778742, 6.929842397649791, 12.45890945479811, 42.89950443515912, 51.97148119386282, 19.91793526767113, 33.0564009201021, 10.2863728020849, 20.09169660146593, 10.07363342512079, 20.61153820561587, 19.81191271120439),
    @(8.442852736007525, 6.928129114549552, 12.46032011997066, 42.90099998429218, 51.96804481541183, 19.92242486565425, 33.05982688327098, 10.28861857022874, 20.05454026875693, 10.07552561032164, 20.59670994498731, 19.8189770162896),
    @(8.426329027936072, 6.918463407308792, 12.46878646623713, 42.9125746019335, 51.95276249246006, 19.94924696072987, 33.08162249122125, 10.30172562318296, 19.84117182399831, 10.0865795217134, 20.51241641754422, 19.86002331455102),
    @(8.416458890730473, 6.912671184455774, 12.47432193060842, 42.92235442788825, 51.94732133661734, 19.96667665663297, 33.0969214703746, 10.30997862042924, 19.70994755527222, 10.09354895656072, 20.46133016754281, 19.88570784403547),
    @(8.410880908192887, 6.90939067037052, 12.47763149656876, 42.92896338819238, 51.94563674465936, 19.97705973483455, 33.10642882383698, 10.31480366538517, 19.63437516686777, 10.09762689161601, 20.43218319129602, 19.90066651633508),
    @(8.409009413043066, 6.908288757347182, 12.47877366542142, 42.93136996959049, 51.94531446108606, 19.98063671041358, 33.10976933452955, 10.31645077652536, 19.60877367106691, 10.09901952523952, 20.42235577729848, 19.90576318348768),
    @(8.417499355398615, 6.913282508111758, 12.47371966378131, 42.92121150447253, 51.94775098241307, 19.96478417699505, 33.09521959635808, 10.30909199082987, 19.72392708279312, 10.09279987894283, 20.46674406226744, 19.88295448128618),
    @(8.447200368914082, 6.930666640380922, 12.45823976084664, 42.89884043989988, 51.97320705290744, 19.91580174983151, 33.05479635994085, 10.28530009085624, 20.10950539801189, 10.07272979098341, 20.61866056612081, 19.8085351406406),
    @(8.467470589677758, 6.942471024754387, 12.4492386985416, 42.89299411148787, 52.00272244566366, 19.88698516411588, 33.03469563338222, 10.27044356812363, 20.3601165979511, 10.06022716976492, 20.71990778092564, 19.76154200473049),
    @(8.456574555186478, 6.936130847112009, 12.45394079878125, 42.89531107523937, 51.9857910420149, 19.90207226138642, 33.04484499436199, 10.27830951957762, 20.22651287008625, 10.06684397506692, 20.66569791246173, 19.78647309496444),
    @(8.417028660528612, 6.91300597435444, 12.47399155250944, 42.92172514762537, 51.94755224287653, 19.96563863786832, 33.09598679660311, 10.30949258594614, 19.71760733665001, 10.09313831538623, 20.4642957436192, 19.88419867746617),
    @(8.37817131041249, 6.890018790032909, 12.5004277199333, 42.98781514750002, 51.96152403888718, 20.04788581591595, 33.1781978464305, 10.34611786801732, 19.16535295209708, 10.12415225416217, 20.25617162530814, 19.99671337016487)
)

$startRow = 2
$startCol = 3   # column C

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
